# Adds 5 newly-collected submissions (rows 116-120) to the raw log sheet
# "八位序列号收集收集结果yd5" (the first worksheet / rId3 / sheet1.xml),
# matching the data appended upstream. The second worksheet
# "八位序列号收集（收集结果）" (the curated/deduplicated results sheet) keeps
# the exact same visible content, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows: A=submitter name, B=submission timestamp (serial date/time,
# same numeric format as the existing rows), C=8-char hex code, D=QQ number.
$newRows = @(
    @{ Row = 116; A = "。";      B = 45983.3904861111; C = "652eb97a"; D = "3435585501" },
    @{ Row = 117; A = "空白";     B = 45983.913900463;  C = "38effb59"; D = "2274948450" },
    @{ Row = 118; A = "唯我所爱"; B = 45983.9352662037; C = "57ded35e"; D = "2583566531" },
    @{ Row = 119; A = "Archive";  B = 45984.485150463;  C = "7eac9117"; D = "1092129794" },
    @{ Row = 120; A = "伪装 🅥";  B = 45984.6307291667; C = "f8b9590f"; D = "85229794" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}
